$wb = $excel.ActiveWorkbook

# --- PCNCode sheet: "Compensation Good" column (C) gets proper-cased,
# non-duplicated display labels instead of re-using the lower-case
# "Title" text. ---
$pcn = $wb.Worksheets.Item("PCNCode")
$pcn.Range("C2").Value = "Papierosy"
$pcn.Range("C3").Value = "PyłTytoiowy"
$pcn.Range("C4").Value = "Kartony"
$pcn.Range("C5").Value = "Krajanka"

# --- Move the active tab / selection from CustomsUnion to PCNCode
# (this also clears the previous tabSelected="1" on CustomsUnion). ---
$pcn.Activate()

# Scroll so row 4 becomes the top visible row (best-effort "topLeftCell").
$excel.ActiveWindow.ScrollRow = 4

$pcn.Range("C9").Select()
